$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 1232
    $ws.Range("F4").Value = 1459
    $ws.Range("F6").Value = 6122
}
